# Apply the scraped-price update to the cryptos list (cells B2:E51).
# Column D ("Price") values that look like plain numbers must be forced to
# text (NumberFormat "@") before assignment so Excel keeps the exact scraped
# string (e.g. "1.000") instead of auto-converting it to a numeric 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.020.68'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.827.93'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.71'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6209'
$ws.Range('E6').Value = '  -6.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.49'
$ws.Range('E8').Value = '  +6.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07483'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2914'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.63'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07614'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').Value = '1.824.62'
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6613'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.94'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009038'
$ws.Range('E17').Value = '  +8.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.954'
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').Value = '29.017.85'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '2.077.39'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '224.09'
$ws.Range('E21').Value = '  -2.02%  '
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.160'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.000'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.40'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.374'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1353'
$ws.Range('E28').Value = '  -3.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.78'
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.492'
$ws.Range('E30').Value = '  -1.28%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.026'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.203'
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.037'
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05213'
$ws.Range('E34').Value = '  -1.97%  '
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('E36').Value = '  +0.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7306'
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').Value = '1.275.12'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01778'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.306'
$ws.Range('E42').Value = '  +6.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8908'
$ws.Range('E43').Value = '  -4.09%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.53'
$ws.Range('D46').Value = '1.975.90'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.32'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.688'
$ws.Range('E50').Value = '  -4.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3953'
$ws.Range('E51').Value = '  -1.72%  '
